$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("mob")
$ws2 = $wb.Worksheets.Item("inf")

# --- sheet "inf": rework the track EAC parameters ---------------------

# New rows 12/13 names first (so shared-string order matches the source edit)
$ws2.Range("A12").Value = "gross_tk_in_hq_track_lifetime"
$ws2.Range("A13").Value = "high_quality_track_price"

# Row 9: crf_track -> interest_rate
$ws2.Range("A9").Value = "interest_rate"
$ws2.Range("C9").Value = "Interest rate used to calculate capital recovery factor (rate)."
$ws2.Range("B9").Value = 0.08

# Row 11: keep name/value, just clarify the description text
$ws2.Range("C11").Value = "Wage cost of maintaining a turnout (USD)."

# New row 12: design tons for a hq track over its lifetime
$ws2.Range("C12").Value = "Design tons for high quality track. Gross tons that a hq track is suposed to support during its lifetime (gross ton-km)."
$ws2.Range("B12").Value = 200000000
$ws2.Range("B12").NumberFormat = "#,##0"
$ws2.Range("B12").HorizontalAlignment = -4108

# New row 13: price of hq track per km
$ws2.Range("C13").Value = "The price of 1km of hight quality track (USD/km)."
$ws2.Range("B13").Value = 800000
$ws2.Range("B13").NumberFormat = "#,##0"
$ws2.Range("B13").HorizontalAlignment = -4108

# Column B widened to fit the new, larger values
$ws2.Columns.Item(2).ColumnWidth = 11.2

# --- sheet views / active tab ------------------------------------------
# "mob" loses its previous selection/scroll state, "inf" becomes the
# active (selected) sheet with A9 selected.
$ws1.Activate()
$ws1.Range("B3").Select()

$ws2.Activate()
$ws2.Range("A9").Select()
